$d = $word.ActiveDocument

$d.Content.Find.Execute("732÷6=122, 0", $true, $false, $false, $false, $false, $true, 1, $false, "679÷8=84, 7", 2) | Out-Null
$d.Content.Find.Execute("925÷8=115, 5", $true, $false, $false, $false, $false, $true, 1, $false, "961÷8=120, 1", 2) | Out-Null
$d.Content.Find.Execute("555÷9=61, 6", $true, $false, $false, $false, $false, $true, 1, $false, "676÷7=96, 4", 2) | Out-Null
$d.Content.Find.Execute("634÷4=158, 2", $true, $false, $false, $false, $false, $true, 1, $false, "536÷6=89, 2", 2) | Out-Null
$d.Content.Find.Execute("380÷8=47, 4", $true, $false, $false, $false, $false, $true, 1, $false, "616÷9=68, 4", 2) | Out-Null
$d.Content.Find.Execute("852÷6=142, 0", $true, $false, $false, $false, $false, $true, 1, $false, "861÷8=107, 5", 2) | Out-Null
$d.Content.Find.Execute("807÷4=201, 3", $true, $false, $false, $false, $false, $true, 1, $false, "469÷4=117, 1", 2) | Out-Null
$d.Content.Find.Execute("178÷7=25, 3", $true, $false, $false, $false, $false, $true, 1, $false, "353÷8=44, 1", 2) | Out-Null
$d.Content.Find.Execute("612÷8=76, 4", $true, $false, $false, $false, $false, $true, 1, $false, "366÷8=45, 6", 2) | Out-Null
$d.Content.Find.Execute("998÷9=110, 8", $true, $false, $false, $false, $false, $true, 1, $false, "616÷6=102, 4", 2) | Out-Null
$d.Content.Find.Execute("568÷8=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "827÷6=137, 5", 2) | Out-Null
$d.Content.Find.Execute("428÷3=142, 2", $true, $false, $false, $false, $false, $true, 1, $false, "466÷4=116, 2", 2) | Out-Null
$d.Content.Find.Execute("643÷5=128, 3", $true, $false, $false, $false, $false, $true, 1, $false, "369÷9=41, 0", 2) | Out-Null
$d.Content.Find.Execute("600÷6=100, 0", $true, $false, $false, $false, $false, $true, 1, $false, "981÷5=196, 1", 2) | Out-Null
$d.Content.Find.Execute("270÷3=90, 0", $true, $false, $false, $false, $false, $true, 1, $false, "484÷4=121, 0", 2) | Out-Null
$d.Content.Find.Execute("848÷9=94, 2", $true, $false, $false, $false, $false, $true, 1, $false, "801÷2=400, 1", 2) | Out-Null
$d.Content.Find.Execute("449÷2=224, 1", $true, $false, $false, $false, $false, $true, 1, $false, "689÷8=86, 1", 2) | Out-Null
$d.Content.Find.Execute("701÷6=116, 5", $true, $false, $false, $false, $false, $true, 1, $false, "679÷3=226, 1", 2) | Out-Null
$d.Content.Find.Execute("323÷4=80, 3", $true, $false, $false, $false, $false, $true, 1, $false, "246÷5=49, 1", 2) | Out-Null
$d.Content.Find.Execute("334÷2=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "359÷3=119, 2", 2) | Out-Null
$d.Content.Find.Execute("482÷5=96, 2", $true, $false, $false, $false, $false, $true, 1, $false, "326÷4=81, 2", 2) | Out-Null
$d.Content.Find.Execute("676÷2=338, 0", $true, $false, $false, $false, $false, $true, 1, $false, "325÷6=54, 1", 2) | Out-Null
$d.Content.Find.Execute("114÷9=12, 6", $true, $false, $false, $false, $false, $true, 1, $false, "183÷5=36, 3", 2) | Out-Null
$d.Content.Find.Execute("506÷5=101, 1", $true, $false, $false, $false, $false, $true, 1, $false, "961÷5=192, 1", 2) | Out-Null
$d.Content.Find.Execute("328÷3=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "637÷8=79, 5", 2) | Out-Null
